$d = $word.ActiveDocument

$pairs = @(
    @("892×3=2676", "242×8=1936"),
    @("984×4=3936", "455×7=3185"),
    @("638×2=1276", "121×9=1089"),
    @("949×3=2847", "303×2=606"),
    @("427×9=3843", "808×2=1616"),
    @("835×2=1670", "567×2=1134"),
    @("228×9=2052", "870×9=7830"),
    @("333×9=2997", "346×2=692"),
    @("679×6=4074", "664×9=5976"),
    @("590×3=1770", "102×4=408"),
    @("902×9=8118", "576×4=2304"),
    @("967×5=4835", "975×6=5850"),
    @("603×9=5427", "843×5=4215"),
    @("160×4=640", "877×3=2631"),
    @("114×8=912", "206×3=618"),
    @("719×8=5752", "846×6=5076"),
    @("749×9=6741", "443×7=3101"),
    @("337×7=2359", "937×2=1874"),
    @("555×3=1665", "104×8=832"),
    @("151×5=755", "423×9=3807"),
    @("298×9=2682", "709×4=2836"),
    @("990×6=5940", "565×8=4520"),
    @("583×2=1166", "428×3=1284"),
    @("925×7=6475", "860×7=6020"),
    @("286×3=858", "522×6=3132")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
